$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data change -----------------------------------------------------------
# Columns Q, R and S (STR_bronzeRewards / STR_silverRewards / STR_goldRewards)
# for every data row (2-31) used to hold "a:b:1,b:c:1,c:d:1" and are now the
# alliance-shrine resource reward string.
$newReward = "resources:wood:10,resources:stone:10,resources:iron:10"
$ws.Range("Q2:S31").Value = $newReward

# --- View / selection changes ----------------------------------------------
# Ruler was previously forced off (showRuler="0"); it is back to the default.
$excel.ActiveWindow.DisplayRuler = $true

# Move the active selection the same way the author's session ended up:
# top-left viewport around column K, active cell Q6.
$ws.Range("K1").Select()
$ws.Range("Q6").Select()
